$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.890.32"
$ws.Range("E2").Value = "  +3.36%  "

$ws.Range("D3").Value = "3.266.63"
$ws.Range("E3").Value = "  +2.46%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.74"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.91"
$ws.Range("E6").Value = "  +5.80%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "3.266.98"
$ws.Range("E9").Value = "  +2.39%  "

$ws.Range("E10").Value = "  +7.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.74"
$ws.Range("E11").Value = "  +2.73%  "

$ws.Range("E12").Value = "  +6.31%  "

$ws.Range("D13").Value = "3.831.67"
$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.45"
$ws.Range("E15").Value = "  +3.96%  "

$ws.Range("D16").Value = "67.815.75"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("E17").Value = "  +3.26%  "

$ws.Range("D18").Value = "3.280.18"
$ws.Range("E18").Value = "  +3.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("E20").Value = "  +5.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.20"
$ws.Range("E21").Value = "  +4.55%  "

$ws.Range("E22").Value = "  +5.12%  "

$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.30"
$ws.Range("E24").Value = "  +3.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.513"
$ws.Range("E25").Value = "  +3.79%  "

$ws.Range("E26").Value = "  +3.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("E28").Value = "  +2.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.67"
$ws.Range("E31").Value = "  +5.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.77"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  +5.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.92"
$ws.Range("E35").Value = "  +4.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  +4.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.82"
$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("E39").Value = "  +3.15%  "

$ws.Range("E40").Value = "  +11.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.76"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.65"
$ws.Range("E42").Value = "  +11.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  +4.84%  "

$ws.Range("D44").Value = "2.694.71"
$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "350.91"
$ws.Range("E45").Value = "  +6.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.73"
$ws.Range("E46").Value = "  +6.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.82"
$ws.Range("E47").Value = "  +3.02%  "

$ws.Range("E48").Value = "  +2.83%  "

$ws.Range("E49").Value = "  +2.46%  "

$ws.Range("E50").Value = "  +5.24%  "

$ws.Range("E51").Value = "  +0.16%  "
